$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title cell
$ws.Range("A1").Value = "Friday Night Magic 2008 (F08)"

# Remove the rest of the card-list rows entirely, but leave row 2 present
# (empty) instead of also deleting it.
$ws.Range("A3:A13").EntireRow.Delete()
$ws.Range("A2").Value = ""
$ws.Range("A2").Style = "Normal"
